$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6644631028175354
$ws.Range("B1").Value = 1.783824324607849
$ws.Range("C1").Value = 2.439590692520142
$ws.Range("D1").Value = 1.734158992767334
$ws.Range("E1").Value = 0.8633520603179932
